# Module 7.pptx edit
# Slide 10, "Text Placeholder 1" shape: the sub-title line that used to read
# "Using Chrome Postman to query an ODATA" is reworded to
# "Using CRMRestBuilder to query an ODATA from Javascript" and split into
# several runs (CRMRestBuilder / Javascript flagged as possible spelling
# errors by the original author's PowerPoint session).

$p = $ppt.ActivePresentation

$slide = $p.Slides.Item(10)
$shape = $slide.Shapes.Item(1)
$textRange = $shape.TextFrame.TextRange

# Replace the whole line with the new wording first (keeps this script
# resilient to the exact length of whatever text is currently there).
$whole = $textRange.Characters(1, $textRange.Text.Length)
$whole.Text = "Using CRMRestBuilder to query an ODATA from Javascript"

# Re-fetch the range/text now that the content has changed, then carve the
# paragraph up into the six runs the final deck has: each
# Characters(start, length) assignment below forces PowerPoint to break a
# new run at that boundary (even though the text itself doesn't change),
# matching the run layout from the diff.
$textRange = $shape.TextFrame.TextRange
$segments = @(
    @{ Start = 7;  Length = 14 },  # "CRMRestBuilder"
    @{ Start = 21; Length = 4 },   # " to "
    @{ Start = 25; Length = 9 },   # "query an "
    @{ Start = 34; Length = 11 },  # "ODATA from "
    @{ Start = 45; Length = 10 }   # "Javascript"
)

foreach ($segment in $segments) {
    $chars = $textRange.Characters($segment.Start, $segment.Length)
    $chars.Text = $chars.Text
}
